# Apply mapping updates to FRFamilyMemberHistoriesLMCDAFHIR.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the Metadata sheet's "Date" value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-22T09:24:45+00:00"

# --- Update the "Mapping Table 1" sheet mapping row for relationship ---
$map1 = $wb.Worksheets.Item("Mapping Table 1")
$map1.Range("A5").Value = "FRCDAAntecedentsFamiliaux.subject"
$map1.Range("D5").Value = "FRFamilyMemberHistoryDocument.relationship"
